# Calculate EIA impervious treated area
# Adds COUNTIF validation (col C) plus mainstem / subshed / watts lookup
# helper columns (D/E/F) to the eia_location_id sheet, and restores the
# view state (selected cells + active sheet) that was captured when the
# workbook was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- eia_location_id (sheet9): new lookup columns C/D/E/F ---
$ws9 = $wb.Worksheets.Item("eia_location_id")

# New header labels for columns D, E, F (shared strings 113-115)
$ws9.Range("D1").Value = 'mainstem_lookup'
$ws9.Range("E1").Value = 'subshed_lookup'
$ws9.Range("F1").Value = 'watts_lookup'

# Column C: COUNTIF of location_id against the location sheet
# Column D/E/F: classify each segment_lookup text as mainstem / subshed / watts
$ws9.Range("C2").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B2)'
$ws9.Range("E2").Value = 'Rock Creek Upper - MS4 - Portal Branch'
$ws9.Range("C3").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B3)'
$ws9.Range("E3").Value = 'Rock Creek Lower - MS4 - Dumbarton Oaks'
$ws9.Range("C4").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B4)'
$ws9.Range("E4").Value = 'Rock Creek Lower - MS4 - Klingle Valley Run'
$ws9.Range("C5").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B5)'
$ws9.Range("E5").Value = 'Rock Creek Lower - MS4 - Melvin Hazen Valley Branch'
$ws9.Range("C6").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B6)'
$ws9.Range("E6").Value = 'Rock Creek Lower - MS4 - Normanstone Creek'
$ws9.Range("C7").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B7)'
$ws9.Range("E7").Value = 'Potomac Upper - MS4 - Battery Kemble Creek'
$ws9.Range("C8").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B8)'
$ws9.Range("E8").Value = 'Anacostia Upper - MS4 - Fort Davis Tributary'
$ws9.Range("C9").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B9)'
$ws9.Range("E9").Value = 'Anacostia Upper - MS4 - Texas Avenue Tributary'
$ws9.Range("C10").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B10)'
$ws9.Range("E10").Value = 'Anacostia Lower - MS4 - Fort Stanton Tributary'
$ws9.Range("C11").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B11)'
$ws9.Range("E11").Value = 'Potomac Upper - MS4 - Dalecarlia Tributary'
$ws9.Range("C12").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B12)'
$ws9.Range("E12").Value = 'Anacostia Upper - MS4 - Pope Branch'
$ws9.Range("C13").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B13)'
$ws9.Range("E13").Value = 'Anacostia Upper - MS4 - Fort Chaplin Tributary'
$ws9.Range("C14").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B14)'
$ws9.Range("E14").Value = 'Anacostia Upper - MS4 - Fort Dupont Tributary'
$ws9.Range("C15").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B15)'
$ws9.Range("E15").Value = 'Rock Creek Upper - MS4 - Soapstone Creek'
$ws9.Range("C16").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B16)'
$ws9.Range("E16").Value = 'Anacostia Upper - MS4 - Nash Run'
$ws9.Range("C17").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B17)'
$ws9.Range("E17").Value = 'Rock Creek Upper - MS4 - Luzon Branch'
$ws9.Range("C18").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B18)'
$ws9.Range("E18").Value = 'Rock Creek Upper - MS4 - Pinehurst Branch'
$ws9.Range("C19").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B19)'
$ws9.Range("E19").Value = 'Rock Creek Upper - MS4 - Fenwick Branch'
$ws9.Range("C20").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B20)'
$ws9.Range("E20").Value = 'Anacostia Upper - MS4 - Hickey Run'
$ws9.Range("C21").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B21)'
$ws9.Range("E21").Value = 'Potomac Upper - MS4 - Foundry Branch'
$ws9.Range("C22").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B22)'
$ws9.Range("E22").Value = 'Rock Creek Upper - MS4 - Broad Branch'
$ws9.Range("C23").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B23)'
$ws9.Range("D23").Value = 'Rock Creek Lower - CSS - Rock Creek'
$ws9.Range("C24").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B24)'
$ws9.Range("F24").Value = 'Anacostia Upper - MS4 - Watts Branch - Upper'
$ws9.Range("C25").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B25)'
$ws9.Range("F25").Value = 'Anacostia Upper - MS4 - Watts Branch - Lower'
$ws9.Range("C26").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B26)'
$ws9.Range("E26").Value = 'Rock Creek Lower - MS4 - Piney Branch'
$ws9.Range("C27").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B27)'
$ws9.Range("E27").Value = 'Potomac Lower - MS4 - Oxon Run'
$ws9.Range("C28").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B28)'
$ws9.Range("D28").Value = 'Rock Creek Upper - MS4 - Rock Creek'
$ws9.Range("C29").Formula = '=COUNTIF(location!$A$2:$A$29,eia_location_id!B29)'
$ws9.Range("D29").Value = 'Rock Creek Lower - MS4 - Rock Creek'

# Approximate column widths for the new columns (engine cannot hit exact AutoFit pixel widths)
$ws9.Columns.Item(3).ColumnWidth = 9.625
$ws9.Columns.Item(4).ColumnWidth = 31.125
$ws9.Columns.Item(5).ColumnWidth = 44.625
$ws9.Columns.Item(6).ColumnWidth = 38.625

# --- Selections / active sheet bookkeeping -------------------------------
# These mirror incidental view-state changes captured in the source diff
# (selected cell moved on a couple of sheets, and the workbook's active
# tab moved from eia_location_id to eia_area when it was last saved).

$wsLocation = $wb.Worksheets.Item("location")
$wsLocation.Range("C28").Select()

$ws9.Range("D33").Select()

$wsEiaArea = $wb.Worksheets.Item("eia_area")
$wsEiaArea.Range("T9").Select()

